# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
# Mirrors the commit: "Created functions to get season record" — the
# existing sheet (player roster, columns A:AC) gets three new trailing
# columns AD:AF with the team's season record repeated on every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastDataRow = 57          # rows 2..57 hold player data (row 1 is the header)
$winsCol   = 30             # AD
$lossesCol = 31             # AE
$tiesCol   = 32             # AF

$wins   = 83
$losses = 79
$ties   = 0

# Seed the new header cells from the existing header style (AC1, s="1")
# so they inherit the same bold/bordered look as the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item(1, $winsCol).Value   = "Wins"
$ws.Cells.Item(1, $lossesCol).Value = "Losses"
$ws.Cells.Item(1, $tiesCol).Value   = "Ties"

for ($row = 2; $row -le $lastDataRow; $row++) {
    $ws.Cells.Item($row, $winsCol).Value   = $wins
    $ws.Cells.Item($row, $lossesCol).Value = $losses
    $ws.Cells.Item($row, $tiesCol).Value   = $ties
}

Write-Output "Added Wins/Losses/Ties columns (AD:AF) for rows 1-$lastDataRow"
